# Add new column 'Servised by' to Card24 sheet
# - Insert header "Servised by" in O1, matching the style used by the other
#   header cells (A1:N1), which is cellXf index 1 (bold, centered, top
#   aligned, thin border) - reuse it by copying formats from N1 instead of
#   re-declaring a brand new style.
# - Also correct N6 from "عطل" to "nan" per the source update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Copy the existing header style (N1) onto the new header cell (O1) so the
# new column heading renders identically to the rest of the header row,
# then set its text.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O1").Value = "Servised by"

# Data rows 2-12 have no "Servised by" value yet (left blank / not yet
# recorded), so nothing further needs to be written for column O there.

# Correction in existing data: N6 changes from "عطل" to "nan"
$ws.Range("N6").Value = "nan"
